$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1 = @(0.6365460719280274, 5.096915102979979, 0.1237739816780862, 0.01896930699493486)
    2 = @(0.9257594155407242, 6.716417910447761, 0.1560913856018062, 0.03299834600909928)
    3 = @(1.222686689997595, 9.473684210526317, 0.2684382142522298, 0.03825234534431345)
    4 = @(0.06895546953592396, 12.82051282051282, 0.1571457784186893, 0.02542506513314081)
    5 = @(0, 0, 0, 0)
    6 = @(0, 0, 0, 0)
    7 = @(0, 0, 0, 0)
    8 = @(2.874967000349976, 11.76470588235294, 0.3729488452576417, 0.03729488452576417)
    9 = @(0, 0, 0.03556762989432234, 0.01778381494716117)
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Cells.Item($row, 4).Value = $rowVals[0]
    $ws.Cells.Item($row, 5).Value = $rowVals[1]
    $ws.Cells.Item($row, 6).Value = $rowVals[2]
    $ws.Cells.Item($row, 7).Value = $rowVals[3]
}
